$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.327.92'
$ws.Range('E2').Value = '  -2.84%  '
$ws.Range('D3').Value = '1.974.29'
$ws.Range('E3').Value = '  -3.35%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.33'
$ws.Range('E5').Value = '  -2.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('E6').Value = '  -5.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.47'
$ws.Range('E7').Value = '  -11.12%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.374'
$ws.Range('E9').Value = '  -6.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '56.15'
$ws.Range('E10').Value = '  -5.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0873'
$ws.Range('E11').Value = '  +9.04%  '
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.850'
$ws.Range('E13').Value = '  -6.48%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.95'
$ws.Range('E14').Value = '  -5.98%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.262.88'
$ws.Range('E15').Value = '  -3.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.69'
$ws.Range('E16').Value = '  -7.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.42'
$ws.Range('E17').Value = '  -5.49%  '
$ws.Range('D18').Value = '1.983.37'
$ws.Range('E18').Value = '  -2.93%  '
$ws.Range('D19').Value = '36.219.56'
$ws.Range('E19').Value = '  -2.77%  '
$ws.Range('D20').Value = '0.0₃0900'
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.31'
$ws.Range('E21').Value = '  -3.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.23'
$ws.Range('E22').Value = '  -5.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.00'
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.50'
$ws.Range('E25').Value = '  -3.39%  '
$ws.Range('E26').Value = '  -2.86%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.75'
$ws.Range('E27').Value = '  -2.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '166.10'
$ws.Range('E28').Value = '  +2.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.94'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('E31').Value = '  -2.48%  '
$ws.Range('E32').Value = '  -1.42%  '
$ws.Range('E33').Value = '  -5.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0645'
$ws.Range('E34').Value = '  +3.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.40'
$ws.Range('E35').Value = '  -5.47%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.05'
$ws.Range('E37').Value = '  -5.54%  '
$ws.Range('E38').Value = '  -1.87%  '
$ws.Range('E39').Value = '  -7.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.92'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0960'
$ws.Range('E41').Value = '  -5.90%  '
$ws.Range('E42').Value = '  -6.59%  '
$ws.Range('E43').Value = '  -4.98%  '
$ws.Range('E44').Value = '  -3.06%  '
$ws.Range('E45').Value = '  -7.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.99'
$ws.Range('E46').Value = '  -8.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.50'
$ws.Range('E47').Value = '  -5.17%  '
$ws.Range('D48').Value = '1.363.12'
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.36'
$ws.Range('E49').Value = '  -5.93%  '
$ws.Range('E50').Value = '  -3.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '45.26'
$ws.Range('E51').Value = '  -3.93%  '
